$d = $word.ActiveDocument
$d.Content.Find.Execute(
    "of experience working with dynamic and collaborative teams ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "of experience working with collaborative teams ",
    2)
